# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted before existing row 112 (which
# pushes the former rows 112-151 down to 113-152, each keeping its original
# data). The freshly inserted row 112 is then populated with the new
# observation for Macroferia Regional de Talca - Mango.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 112, shifting 112:151 down
# to 113:152.
$ws.Rows.Item(112).Insert()

# Fill in the new row 112 with the new data point.
$ws.Range("A112").Value = 5
$ws.Range("B112").Value = "Macroferia Regional de Talca"
$ws.Range("C112").Value = "Maule"
$ws.Range("D112").Value = 44841
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100108
$ws.Range("H112").Value = "Tropicales y subtropicales"
$ws.Range("I112").Value = 100108002
$ws.Range("J112").Value = "Mango"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 240
$ws.Range("N112").Value = 7000
$ws.Range("O112").Value = 7000
$ws.Range("P112").Value = 7000
$ws.Range("Q112").Value = "`$/bandeja 4 kilos"
$ws.Range("R112").Value = "Brasil"
$ws.Range("S112").Value = 1750
$ws.Range("T112").Value = 4
